$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 already carries the "header" style (s=1) in the source workbook, so copy
# its formatting into A1 before the rest of the layout changes, rather than
# rebuilding the bold/border/alignment combo from scratch (which would create
# a brand-new style entry instead of reusing the existing one).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Clear the old D column / extra empty rows so the sheet matches the new
# trimmed A1:C4 layout.
$ws.Range("D1:D5").Clear()
$ws.Range("A4:A5").Clear()

# A2/A3 previously carried the header style (s=1) too; the new layout keeps
# that styling only on the header row, so strip formatting from the old
# counter cells (new counter values go in unstyled cells).
$ws.Range("A2:A3").ClearFormats()

# Header row values
$ws.Range("A1").Value = "counter"
$ws.Range("B1").Value = "participant 1"
$ws.Range("C1").Value = "participant 2"

# Data rows (counter now starts at 1 instead of 0)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Oxy"
$ws.Range("C2").Value = "Liu"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Liu"
$ws.Range("C3").Value = "Law"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Angie"
$ws.Range("C4").Value = "Teddy"

$ws.Range("C4").Select()
